# Spilamberto.xlsx - "aggiornamento fino al 26/03"
# Append 5 new daily rows (234-238) below the existing data (which ends at row 233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the formatting of the last existing data row (233) onto the
# new rows (234-238) so the new cells inherit the same style (bold,
# bordered, centered date format for column A) as the rest of the table.
$ws.Range("A233:D233").Copy()
$ws.Range("A234:D238").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data: date (Excel serial), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(234, 44308, 3, 23, 180.6046329014527),
    @(235, 44309, 1, 22, 172.7522575579113),
    @(236, 44310, 1, 18, 141.3427561837456),
    @(237, 44311, 5, 21, 164.8998822143699),
    @(238, 44312, 5, 22, 172.7522575579113)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}
